# data for csc computations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data row: graph pointer
$ws.Range("A9").Value = "graph"
$ws.Range("B9").Value = "data/canal_network_matrix_50meters.p"

# Update selection to E11
$ws.Range("E11").Select()

# Update window position
$excel.ActiveWindow.Left = 2175
$excel.ActiveWindow.Top = 1395
